$d = $word.ActiveDocument

# Locate the "Draft 1: Concise IEEE Style" heading - the whole first
# draft (and everything up to and including the "Draft 2: Academic
# Descriptive Style" heading that introduces the kept draft) is removed,
# leaving the reviewed "4.3 Software Interfaces" section from Draft 2 in
# place as the sole remaining content.
$startRange = $d.Content
$startRange.Find.Execute("Draft 1: Concise IEEE Style", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $startRange.Start

$endRange = $d.Content
$endRange.Find.Execute("Draft 2: Academic Descriptive Style", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
# Extend to the end of that heading's paragraph (consumes its paragraph
# mark too) so the following "4.3 Software Interfaces" paragraph becomes
# the first paragraph in the document.
$endPara = $d.Range($endRange.Start, $endRange.Start).Paragraphs(1)
$endPos = $endPara.Range.End

$rng = $d.Range($startPos, $endPos)
$rng.Delete()
